$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D/E columns to keep text type while we write values that look numeric
# (Excel auto-converts numeric-looking strings assigned to .Value unless the
# cell is already formatted as Text). We restore formatting afterwards so the
# cell style stays the same as before (plain, unstyled data cells).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.945.01"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "1.919.19"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "239.64"
$ws.Range("E5").Value = "  -3.37%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.4918"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "0.2967"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.06780"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "1.909.27"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").Value = "17.03"
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").Value = "0.07316"
$ws.Range("D13").Value = "5.149"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "90.07"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "0.6733"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "30.908.50"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "0.000007953"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "13.45"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "2.168.27"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "5.186"
$ws.Range("E22").Value = "  +6.64%  "
$ws.Range("D23").Value = "207.86"
$ws.Range("E23").Value = "  +7.33%  "
$ws.Range("D24").Value = "6.244"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").Value = "9.687"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").Value = "158.12"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").Value = "18.89"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").Value = "1.981"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("D29").Value = "1.426"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").Value = "4.324"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "0.09178"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").Value = "4.076"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "0.05179"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").Value = "0.7532"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").Value = "1.121"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").Value = "2.736"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "0.01853"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").Value = "2.737"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").Value = "0.9250"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").Value = "2.102"
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("D41").Value = "0.4533"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("D42").Value = "107.68"
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").Value = "5.898"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("D44").Value = "1.010"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").Value = "0.1409"
$ws.Range("E45").Value = "  +4.56%  "
$ws.Range("D46").Value = "7.727"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Value = "66.63"
$ws.Range("E47").Value = "  +14.51%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "35.14"
$ws.Range("E48").Value = "  +4.44%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "0.4101"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05952"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "8.941"
$ws.Range("E51").Value = "  +1.80%  "

$dataRange.ClearFormats()
